$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2039.375
$ws.Range("J17").Value = 2039.375
$ws.Range("L17").Value = 6118.125
$ws.Range("N17").Value = -6454.125
$ws.Range("H43").Value = 1894.7273
$ws.Range("I43").Value = 838.5
$ws.Range("K43").Value = 838.5
$ws.Range("M43").Value = -769.5
$ws.Range("H96").Value = 1128.0714
$ws.Range("I96").Value = 397.08334
$ws.Range("J96").Value = 5514
$ws.Range("K96").Value = 1191.25002
$ws.Range("L96").Value = 16542
$ws.Range("M96").Value = 181.7499800000001
$ws.Range("N96").Value = -19288
$ws.Range("H97").Value = 2263.4
$ws.Range("J97").Value = 2402
$ws.Range("L97").Value = 7206
$ws.Range("N97").Value = -8198
$ws.Range("H132").Value = 251854.61
$ws.Range("I132").Value = 298919.34
$ws.Range("J132").Value = 16531
$ws.Range("K132").Value = 896758.02
$ws.Range("L132").Value = 49593
$ws.Range("M132").Value = -894228.02
$ws.Range("N132").Value = -54653
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 759.5599999999999
$ws.Range("I2").Value = 639.8158
$ws.Range("K2").Value = 639.8158
$ws.Range("M2").Value = -526.8158
$ws.Range("H32").Value = 2099.97
$ws.Range("I32").Value = 2099.97
$ws.Range("K32").Value = 2099.97
$ws.Range("M32").Value = -1812.97
$ws.Range("H61").Value = 3924.9062
$ws.Range("I61").Value = 3380.3462
$ws.Range("K61").Value = 3380.3462
$ws.Range("M61").Value = -3168.3462
$ws.Range("H74").Value = 20834696
$ws.Range("I74").Value = 25000936
$ws.Range("K74").Value = 25000936
$ws.Range("M74").Value = -25000062
$ws.Range("H77").Value = 20834696
$ws.Range("I77").Value = 25000936
$ws.Range("K77").Value = 125004680
$ws.Range("M77").Value = -125000312
$ws.Range("H88").Value = 3522.5715
$ws.Range("I88").Value = 1498.5
$ws.Range("K88").Value = 1498.5
$ws.Range("M88").Value = -1092.5
$ws.Range("H91").Value = 3522.5715
$ws.Range("I91").Value = 1498.5
$ws.Range("K91").Value = 1498.5
$ws.Range("M91").Value = -94.5
$ws.Range("H97").Value = 538.17645
$ws.Range("I97").Value = 563.38464
$ws.Range("J97").Value = 456.25
$ws.Range("K97").Value = 563.38464
$ws.Range("L97").Value = 456.25
$ws.Range("M97").Value = -67.38463999999999
$ws.Range("N97").Value = -1448.25
$ws.Range("H102").Value = 1507.95
$ws.Range("I102").Value = 1374.7059
$ws.Range("J102").Value = 2263
$ws.Range("K102").Value = 1374.7059
$ws.Range("L102").Value = 2263
$ws.Range("M102").Value = 247.2941000000001
$ws.Range("N102").Value = -5507
$ws.Range("H116").Value = 759.5599999999999
$ws.Range("I116").Value = 639.8158
$ws.Range("K116").Value = 639.8158
$ws.Range("M116").Value = 1654.1842
$ws.Range("H122").Value = 3660.65
$ws.Range("I122").Value = 2741.2666
$ws.Range("K122").Value = 8223.799800000001
$ws.Range("M122").Value = -5773.799800000001
$ws.Range("H132").Value = 17061.805
$ws.Range("I132").Value = 19375.94
$ws.Range("K132").Value = 58127.81999999999
$ws.Range("M132").Value = -55597.81999999999
$ws.Range("H136").Value = 3924.9062
$ws.Range("I136").Value = 3380.3462
$ws.Range("K136").Value = 10141.0386
$ws.Range("M136").Value = -7591.0386
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 759.5599999999999
$ws.Range("I3").Value = 639.8158
$ws.Range("K3").Value = 639.8158
$ws.Range("M3").Value = -525.8158
$ws.Range("H94").Value = 1307.04
$ws.Range("I94").Value = 1428.7778
$ws.Range("J94").Value = 994
$ws.Range("K94").Value = 1428.7778
$ws.Range("L94").Value = 994
$ws.Range("M94").Value = -977.7778000000001
$ws.Range("N94").Value = -1896
$ws.Range("H134").Value = 2455.4888
$ws.Range("I134").Value = 2428.6191
$ws.Range("J134").Value = 2831.6667
$ws.Range("K134").Value = 7285.8573
$ws.Range("L134").Value = 8495.000100000001
$ws.Range("M134").Value = -4750.8573
$ws.Range("N134").Value = -13565.0001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7279.7
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 7279.7
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 7279.7
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -7869.7
$ws.Range("H34").Value = 7279.7
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 7279.7
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 7279.7
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -7683.7
$ws.Range("H97").Value = 6222.222
$ws.Range("J97").Value = 6625
$ws.Range("L97").Value = 6625
$ws.Range("N97").Value = -8607
$ws.Range("H132").Value = 32522478
$ws.Range("I132").Value = 37038870
$ws.Range("J132").Value = 4439
$ws.Range("K132").Value = 111116610
$ws.Range("L132").Value = 13317
$ws.Range("M132").Value = -111114080
$ws.Range("N132").Value = -18377
$ws.Range("H134").Value = 1942.2565
$ws.Range("I134").Value = 1914.421
$ws.Range("K134").Value = 5743.263
$ws.Range("M134").Value = -3208.263
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1114.8
$ws.Range("H132").Value = 2559
$ws.Range("I132").Value = 1241.4445
$ws.Range("K132").Value = 11173.0005
$ws.Range("M132").Value = -8643.0005
$ws.Range("H141").Value = 6151.6523
$ws.Range("I141").Value = 3584.5386
$ws.Range("J141").Value = 9488.9
$ws.Range("K141").Value = 10753.6158
$ws.Range("L141").Value = 28466.7
$ws.Range("M141").Value = -5573.6158
$ws.Range("N141").Value = -38826.7
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1488.875
$ws.Range("I97").Value = 1260.3529
$ws.Range("J97").Value = 2043.8572
$ws.Range("K97").Value = 1260.3529
$ws.Range("L97").Value = 2043.8572
$ws.Range("M97").Value = -764.3529000000001
$ws.Range("N97").Value = -3035.8572
$ws.Range("H123").Value = 49011.848
$ws.Range("J123").Value = 49011.848
$ws.Range("L123").Value = 49011.848
$ws.Range("N123").Value = -53911.848
$ws.Range("H132").Value = 145171.5
$ws.Range("I132").Value = 222934.11
$ws.Range("K132").Value = 668802.33
$ws.Range("M132").Value = -666272.33
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 12581.5
$ws.Range("J69").Value = 12581.5
$ws.Range("L69").Value = 12581.5
$ws.Range("N69").Value = -14203.5
$ws.Range("H72").Value = 12581.5
$ws.Range("J72").Value = 12581.5
$ws.Range("L72").Value = 37744.5
$ws.Range("N72").Value = -45856.5
$ws.Range("H74").Value = 30000
$ws.Range("I74").Value = 30000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 30000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -29002
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 30000
$ws.Range("I77").Value = 30000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 90000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -85008
$ws.Range("N77").Value = ""
$ws.Range("H96").Value = 17499.5
$ws.Range("J96").Value = 17499.5
$ws.Range("L96").Value = 17499.5
$ws.Range("N96").Value = -22991.5
$ws.Range("H132").Value = 5480.3286
$ws.Range("I132").Value = 4730.4526
$ws.Range("K132").Value = 14191.3578
$ws.Range("M132").Value = -11661.3578
$ws.Range("H136").Value = 3583.25
$ws.Range("I136").Value = 2529.0625
$ws.Range("K136").Value = 7587.1875
$ws.Range("M136").Value = -5037.1875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5482.5
$ws.Range("J14").Value = 5482.5
$ws.Range("L14").Value = 5482.5
$ws.Range("N14").Value = -5818.5
$ws.Range("H54").Value = 38594.6
$ws.Range("J54").Value = 38243.25
$ws.Range("L54").Value = 38243.25
$ws.Range("N54").Value = -39283.25
$ws.Range("H76").Value = 39999
$ws.Range("J76").Value = 39999
$ws.Range("L76").Value = 39999
$ws.Range("N76").Value = -40629
$ws.Range("H79").Value = 39999
$ws.Range("J79").Value = 39999
$ws.Range("L79").Value = 39999
$ws.Range("N79").Value = -42183
$ws.Range("H96").Value = 1959.9231
$ws.Range("J96").Value = 2959.2
$ws.Range("L96").Value = 2959.2
$ws.Range("N96").Value = -5705.2
$ws.Range("H100").Value = 10490.637
$ws.Range("I100").Value = 850
$ws.Range("K100").Value = 1700
$ws.Range("M100").Value = -1159
$ws.Range("H132").Value = 11497494
$ws.Range("I132").Value = 23811980
$ws.Range("J132").Value = 3972.8
$ws.Range("K132").Value = 71435940
$ws.Range("L132").Value = 11918.4
$ws.Range("M132").Value = -71433410
$ws.Range("N132").Value = -16978.4
